# Adds the 17 new MIMARKS survey/built attribute columns (AC..AS) to row 15,
# along with their header text and explanatory cell comments, matching the
# upstream "harmonized name" field additions.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newColumns = @(
    @{ Col = "AC"; Header = "dew_point"; Comment = "temperature to which a given parcel of humid air must be cooled, at constant barometric pressure, for water vapor to condense into water." },
    @{ Col = "AD"; Header = "indoor_surf"; Comment = "type of indoor surface" },
    @{ Col = "AE"; Header = "isolation_source"; Comment = "Describes the physical, environmental and/or local geographical source of the biological sample from which the sample was derived." },
    @{ Col = "AF"; Header = "rel_to_oxygen"; Comment = "Aerobic or anaerobic" },
    @{ Col = "AG"; Header = "samp_collect_device"; Comment = "Method or device employed for collecting sample" },
    @{ Col = "AH"; Header = "samp_mat_process"; Comment = "Processing applied to the sample during or after isolation" },
    @{ Col = "AI"; Header = "samp_size"; Comment = "Amount or size of sample (volume, mass or area) that was collected" },
    @{ Col = "AJ"; Header = "samp_sort_meth"; Comment = "method by which samples are sorted" },
    @{ Col = "AK"; Header = "samp_vol_we_dna_ext"; Comment = "volume (mL) or weight (g) of sample processed for DNA extraction" },
    @{ Col = "AL"; Header = "source_material_id"; Comment = "unique identifier assigned to a material sample used for extracting nucleic acids, and subsequent sequencing. The identifier can refer either to the original material collected or to any derived sub-samples." },
    @{ Col = "AM"; Header = "substructure_type"; Comment = "substructure or under building is that largely hidden section of the building which is built off the foundations to the ground floor level" },
    @{ Col = "AN"; Header = "surf_air_cont"; Comment = "contaminant identified on surface" },
    @{ Col = "AO"; Header = "surf_humidity"; Comment = "surfaces: water activity as a function of air and material moisture" },
    @{ Col = "AP"; Header = "surf_material"; Comment = "surface materials at the point of sampling" },
    @{ Col = "AQ"; Header = "surf_moisture"; Comment = "water held on a surface" },
    @{ Col = "AR"; Header = "surf_moisture_ph"; Comment = "pH measurement of surface" },
    @{ Col = "AS"; Header = "surf_temp"; Comment = "temperature of the surface at the time of sampling" }
)

# Use an existing optional ("yellow") header cell as the formatting template
# so the new header cells pick up the same font/fill/border/alignment style.
$templateCell = $ws.Range("C15")
$templateCell.Copy()

foreach ($item in $newColumns) {
    $cellRef = $item.Col + "15"
    $cell = $ws.Range($cellRef)

    $cell.PasteSpecial(-4122)
    $cell.Value = $item.Header
    $cell.AddComment($item.Comment)
}

$excel.CutCopyMode = $false

Write-Output "Added $($newColumns.Count) new columns (AC15:AS15) with headers and comments"
